# Add numero contrat generation function
# Regenerates the "etat des taxes" rows: updates the existing tenant rows,
# inserts a new "Logement de fonction" row for the same tiers, and keeps the
# trailing totals row in sync one row further down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper so each contract number is produced the same way (numero contrat
# generation): "<seq>/TEST DR" for the regional-direction rows and
# "<seq>/LF/TEST DR" for the linked logement-de-fonction row.
function New-NumeroContrat([string]$seq, [string]$suffix) {
    if ($suffix) {
        return "$seq/$suffix/TEST DR"
    }
    return "$seq/TEST DR"
}

$seq = "001"

# Insert a fresh row above the old totals row (row 5) so the totals row
# shifts from row 5 down to row 7, leaving room for the new
# "Logement de fonction" entry at row 6.
$ws.Rows("6:6").Insert()

# --- Row 2 : existing "Direction régionale" tenant, updated amounts ---
$ws.Range("A2").Value = (New-NumeroContrat $seq $null)
$ws.Range("B2").Value = "Direction régionale"
$ws.Range("C2").Value = "BB779645"
$ws.Range("D2").Value = "Karami abdelilah"
$ws.Range("E2").Value = "non"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = "--"
$ws.Range("I2").Value = 5000
$ws.Range("J2").Value = "--"
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = 9500

# --- Row 3 : existing "Direction régionale" tenant, updated amounts ---
$ws.Range("A3").Value = (New-NumeroContrat $seq $null)
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "BG1949"
$ws.Range("D3").Value = "Ahmed Test"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = "--"
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = "--"
$ws.Range("K3").Value = 300
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 5700

# --- Row 4 : now also a "Direction régionale" tenant (was "Siège") ---
$ws.Range("A4").Value = (New-NumeroContrat $seq $null)
$ws.Range("B4").Value = "Direction régionale"
# leading apostrophe forces the numeric-looking CIN/IF to be stored as text
# (matches column C's existing text-typed entries), then the style is reset
# so no stray "quote prefix" formatting sticks to the cell.
$ws.Range("C4").Value = "'1196797"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "mediexpets"
$ws.Range("E4").Value = "oui"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = "--"
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = 2000

# --- Row 5 : new "Direction régionale" tenant (replaces the blank spacer) ---
$ws.Range("A5").Value = (New-NumeroContrat $seq $null)
$ws.Range("B5").Value = "Direction régionale"
$ws.Range("C5").Value = "BJ49785"
$ws.Range("D5").Value = "Anas tawfiqi"
$ws.Range("E5").Value = "non"
$ws.Range("F5").Value = "mensuelle"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = "--"
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = 2000

# --- Row 6 : new "Logement de fonction" row, linked to the same tiers as row 4 ---
$ws.Range("A6").Value = (New-NumeroContrat $seq "LF")
$ws.Range("B6").Value = "Logement de fonction"
$ws.Range("C6").Value = "'1196797"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "mediexpets"
$ws.Range("E6").Value = "oui"
$ws.Range("F6").Value = "mensuelle"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = "--"
$ws.Range("I6").Value = 90000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 120000

# --- Row 7 : totals row (was row 5), recalculated for the new rows ---
$ws.Range("A7").Value = " "
$ws.Range("B7").Value = " "
$ws.Range("C7").Value = " "
$ws.Range("D7").Value = " "
$ws.Range("E7").Value = " "
$ws.Range("F7").Value = " "
$ws.Range("G7").Value = " "
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 100000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 800
$ws.Range("L7").Value = 40000
$ws.Range("M7").Value = 139200
